$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prefix table (rows 3-9, 11-12): move the description text
#     from column D into column E, leaving D empty. ---
$prefixRows = 3,4,5,6,7,8,9,11,12
foreach ($r in $prefixRows) {
    $dCell = $ws.Cells.Item($r, 4)
    $text = $dCell.Text
    $ws.Cells.Item($r, 5).Value = $text
    $dCell.Value = $null
}

# --- Row 21: refreshed "last modified" timestamp ---
$ws.Range("B21").Value = "2023-08-17T12:23:53+00:00"

# --- Row 23 header: re-order D/E/F so "skos:broader" moves in front
#     of "skos:definition"/"dct:source" (shifted one column right). ---
$ws.Range("D23").Value = 'skos:broader(separator=",")'
$ws.Range("E23").Value = "skos:definition@en"
$ws.Range("F23").Value = 'dct:source(separator=",")'

# --- Data rows: the "broader concept" values that used to live in
#     column F actually belong in column D; clear the old F value. ---
$broaderMap = @{
    24 = "vocab:1000"
    25 = "vocab:1000"
    26 = "vocab:1000"
    27 = "vocab:1000,vocab:1007"
    28 = "vocab:1000,vocab:1007"
    29 = "vocab:1000,vocab:1007"
    30 = "vocab:1000,vocab:1007"
    31 = "vocab:1007"
    32 = "vocab:1007"
    33 = "vocab:1007"
    34 = "vocab:1007"
    35 = "vocab:1007"
    36 = "vocab:1007"
    37 = "vocab:1007"
}

foreach ($r in $broaderMap.Keys) {
    $ws.Cells.Item($r, 4).Value = $broaderMap[$r]
    $ws.Cells.Item($r, 6).Value = $null
}
